$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 343.9479473333333
$ws.Range("H2").Value = 1031.843842
$ws.Range("I2").Value = 0.9666099193889262
$ws.Range("J2").Value = 0.966609919388926
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 50.862619
$ws.Range("N2").Value = 152.587857
$ws.Range("O2").Value = 0.6466984659960481
$ws.Range("P2").Value = 0.646698465996048
$ws.Range("Q2").Value = 17494.0934010474
$ws.Range("R2").Value = 157446.8406094266
$ws.Range("S2").Value = 0.6251051520853822
$ws.Range("T2").Value = 0.6251051520853821
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 343.9479473333333
$ws.Range("H3").Value = 1031.843842
$ws.Range("I3").Value = 0.9666099193889262
$ws.Range("J3").Value = 0.966609919388926
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.214243
$ws.Range("N3").Value = 12.642729
$ws.Range("O3").Value = 0.05358246462759977
$ws.Range("P3").Value = 0.05358246462759976
$ws.Range("Q3").Value = 1449.480229413869
$ws.Range("R3").Value = 13045.32206472482
$ws.Range("S3").Value = 0.0517933418143442
$ws.Range("T3").Value = 0.05179334181434418
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 343.9479473333333
$ws.Range("H4").Value = 1031.843842
$ws.Range("I4").Value = 0.9666099193889262
$ws.Range("J4").Value = 0.966609919388926
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.01971966666666
$ws.Range("N4").Value = 69.05915899999999
$ws.Range("O4").Value = 0.2926868039589623
$ws.Range("P4").Value = 0.2926868039589623
$ws.Range("Q4").Value = 7917.585327538764
$ws.Range("R4").Value = 71258.26794784887
$ws.Range("S4").Value = 0.282913967980975
$ws.Range("T4").Value = 0.282913967980975
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 343.9479473333333
$ws.Range("H5").Value = 1031.843842
$ws.Range("I5").Value = 0.9666099193889262
$ws.Range("J5").Value = 0.966609919388926
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.5530853333333333
$ws.Range("N5").Value = 1.659256
$ws.Range("O5").Value = 0.007032265417389923
$ws.Range("P5").Value = 0.007032265417389922
$ws.Range("Q5").Value = 190.2325651001725
$ws.Range("R5").Value = 1712.093085901552
$ws.Range("S5").Value = 0.006797457508224807
$ws.Range("T5").Value = 0.006797457508224805
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.970184
$ws.Range("H6").Value = 17.910552
$ws.Range("I6").Value = 0.01677823379880302
$ws.Range("J6").Value = 0.01677823379880302
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 50.862619
$ws.Range("N6").Value = 152.587857
$ws.Range("O6").Value = 0.6466984659960481
$ws.Range("P6").Value = 0.646698465996048
$ws.Range("Q6").Value = 303.659194151896
$ws.Range("R6").Value = 2732.932747367064
$ws.Range("S6").Value = 0.01085045805980896
$ws.Range("T6").Value = 0.01085045805980896
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.970184
$ws.Range("H7").Value = 17.910552
$ws.Range("I7").Value = 0.01677823379880302
$ws.Range("J7").Value = 0.01677823379880302
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.214243
$ws.Range("N7").Value = 12.642729
$ws.Range("O7").Value = 0.05358246462759977
$ws.Range("P7").Value = 0.05358246462759976
$ws.Range("Q7").Value = 25.159806130712
$ws.Range("R7").Value = 226.438255176408
$ws.Range("S7").Value = 0.0008990191190379619
$ws.Range("T7").Value = 0.0008990191190379618
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.970184
$ws.Range("H8").Value = 17.910552
$ws.Range("I8").Value = 0.01677823379880302
$ws.Range("J8").Value = 0.01677823379880302
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 23.01971966666666
$ws.Range("N8").Value = 69.05915899999999
$ws.Range("O8").Value = 0.2926868039589623
$ws.Range("P8").Value = 0.2926868039589623
$ws.Range("Q8").Value = 137.4319620384186
$ws.Range("R8").Value = 1236.887658345768
$ws.Range("S8").Value = 0.004910767626647897
$ws.Range("T8").Value = 0.004910767626647896
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.970184
$ws.Range("H9").Value = 17.910552
$ws.Range("I9").Value = 0.01677823379880302
$ws.Range("J9").Value = 0.01677823379880302
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5530853333333333
$ws.Range("N9").Value = 1.659256
$ws.Range("O9").Value = 0.007032265417389923
$ws.Range("P9").Value = 0.007032265417389922
$ws.Range("Q9").Value = 3.302021207701333
$ws.Range("R9").Value = 29.718190869312
$ws.Range("S9").Value = 0.0001179889933082053
$ws.Range("T9").Value = 0.0001179889933082052
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.671367333333334
$ws.Range("H10").Value = 8.014102000000001
$ws.Range("I10").Value = 0.007507444608265281
$ws.Range("J10").Value = 0.00750744460826528
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 50.862619
$ws.Range("N10").Value = 152.587857
$ws.Range("O10").Value = 0.6466984659960481
$ws.Range("P10").Value = 0.646698465996048
$ws.Range("Q10").Value = 135.8727388843793
$ws.Range("R10").Value = 1222.854649959414
$ws.Range("S10").Value = 0.00485505291171546
$ws.Range("T10").Value = 0.004855052911715458
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.671367333333334
$ws.Range("H11").Value = 8.014102000000001
$ws.Range("I11").Value = 0.007507444608265281
$ws.Range("J11").Value = 0.00750744460826528
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.214243
$ws.Range("N11").Value = 12.642729
$ws.Range("O11").Value = 0.05358246462759977
$ws.Range("P11").Value = 0.05358246462759976
$ws.Range("Q11").Value = 11.25779108492867
$ws.Range("R11").Value = 101.320119764358
$ws.Range("S11").Value = 0.000402267385166039
$ws.Range("T11").Value = 0.0004022673851660389
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.671367333333334
$ws.Range("H12").Value = 8.014102000000001
$ws.Range("I12").Value = 0.007507444608265281
$ws.Range("J12").Value = 0.00750744460826528
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 23.01971966666666
$ws.Range("N12").Value = 69.05915899999999
$ws.Range("O12").Value = 0.2926868039589623
$ws.Range("P12").Value = 0.2926868039589623
$ws.Range("Q12").Value = 61.49412714002423
$ws.Range("R12").Value = 553.447144260218
$ws.Range("S12").Value = 0.002197329968292109
$ws.Range("T12").Value = 0.002197329968292108
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.671367333333334
$ws.Range("H13").Value = 8.014102000000001
$ws.Range("I13").Value = 0.007507444608265281
$ws.Range("J13").Value = 0.00750744460826528
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.5530853333333333
$ws.Range("N13").Value = 1.659256
$ws.Range("O13").Value = 0.007032265417389923
$ws.Range("P13").Value = 0.007032265417389922
$ws.Range("Q13").Value = 1.477494092012445
$ws.Range("R13").Value = 13.297446828112
$ws.Range("S13").Value = 0.00005279434309167438
$ws.Range("T13").Value = 0.00005279434309167436
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.239611333333334
$ws.Range("H14").Value = 9.718834000000001
$ws.Range("I14").Value = 0.009104402204005551
$ws.Range("J14").Value = 0.00910440220400555
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 50.862619
$ws.Range("N14").Value = 152.587857
$ws.Range("O14").Value = 0.6466984659960481
$ws.Range("P14").Value = 0.646698465996048
$ws.Range("Q14").Value = 164.7751169554153
$ws.Range("R14").Value = 1482.976052598738
$ws.Range("S14").Value = 0.005887802939141429
$ws.Range("T14").Value = 0.005887802939141427
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.239611333333334
$ws.Range("H15").Value = 9.718834000000001
$ws.Range("I15").Value = 0.009104402204005551
$ws.Range("J15").Value = 0.00910440220400555
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.214243
$ws.Range("N15").Value = 12.642729
$ws.Range("O15").Value = 0.05358246462759977
$ws.Range("P15").Value = 0.05358246462759976
$ws.Range("Q15").Value = 13.65250938422067
$ws.Range("R15").Value = 122.872584457986
$ws.Range("S15").Value = 0.0004878363090515688
$ws.Range("T15").Value = 0.0004878363090515686
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.239611333333334
$ws.Range("H16").Value = 9.718834000000001
$ws.Range("I16").Value = 0.009104402204005551
$ws.Range("J16").Value = 0.00910440220400555
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 23.01971966666666
$ws.Range("N16").Value = 69.05915899999999
$ws.Range("O16").Value = 0.2926868039589623
$ws.Range("P16").Value = 0.2926868039589623
$ws.Range("Q16").Value = 74.57494472228956
$ws.Range("R16").Value = 671.174502500606
$ws.Range("S16").Value = 0.002664738383047318
$ws.Range("T16").Value = 0.002664738383047316
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.239611333333334
$ws.Range("H17").Value = 9.718834000000001
$ws.Range("I17").Value = 0.009104402204005551
$ws.Range("J17").Value = 0.00910440220400555
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.5530853333333333
$ws.Range("N17").Value = 1.659256
$ws.Range("O17").Value = 0.007032265417389923
$ws.Range("P17").Value = 0.007032265417389922
$ws.Range("Q17").Value = 1.791781514167111
$ws.Range("R17").Value = 16.126033627504
$ws.Range("S17").Value = 0.00006402457276523683
$ws.Range("T17").Value = 0.00006402457276523682
